$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Account Name" column (C) with values
$ws.Range("C1").Value = "Account Name"
$ws.Range("C2").Value = "Account1"
$ws.Range("C3").Value = "Account2"
$ws.Range("C4").Value = "Account3"
$ws.Range("C5").Value = "Account1"

# Autofit the new column width to match bestFit behavior
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Update the active selection to reflect the new last column (C6)
$ws.Range("C6").Select() | Out-Null
